# Updates cryptos list values (price + 1h volume change) to match latest
# snapshot, and swaps the THORChain/PaxDollar row order (ranks 48/49).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" cells are numeric-looking strings (e.g. "232.60", "0.620")
# that would otherwise be auto-coerced to numbers (losing the exact
# formatting, e.g. trailing zeros) when assigned via .Value. Force those
# specific cells to Text format first so the literal string is preserved,
# matching the source data which stores them as plain text.
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D13",
    "D16",
    "D18",
    "D20",
    "D21",
    "D25",
    "D26",
    "D27",
    "D32",
    "D33",
    "D36",
    "D39",
    "D42",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "34.947.67"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.843.28"
$ws.Range("E3").Value = "  +2.19%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "232.60"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "0.620"
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "41.21"
$ws.Range("E8").Value = "  +6.49%  "
$ws.Range("E9").Value = "  +4.09%  "
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "2.110.06"
$ws.Range("E12").Value = "  +2.11%  "
$ws.Range("D13").Value = "11.37"
$ws.Range("E13").Value = "  +5.13%  "
$ws.Range("D14").Value = "1.837.81"
$ws.Range("E14").Value = "  +1.91%  "
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("D16").Value = "4.66"
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").Value = "34.968.05"
$ws.Range("D18").Value = "69.90"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "240.35"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").Value = "12.15"
$ws.Range("E21").Value = "  +4.34%  "
$ws.Range("E22").Value = "  +3.34%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("D25").Value = "171.78"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").Value = "7.82"
$ws.Range("E26").Value = "  +1.60%  "
$ws.Range("D27").Value = "17.45"
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("E28").Value = "  +4.54%  "
$ws.Range("E29").Value = "  +8.98%  "
$ws.Range("E30").Value = "  +0.01%  "
$ws.Range("E31").Value = "  +2.28%  "
$ws.Range("D32").Value = "3.96"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "3.89"
$ws.Range("E33").Value = "  -0.40%  "
$ws.Range("E34").Value = "  +22.69%  "
$ws.Range("E35").Value = "  +11.16%  "
$ws.Range("D36").Value = "0.742"
$ws.Range("E36").Value = "  +9.03%  "
$ws.Range("E37").Value = "  +5.14%  "
$ws.Range("E38").Value = "  +11.91%  "
$ws.Range("D39").Value = "89.84"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").Value = "1.346.96"
$ws.Range("E40").Value = "  +3.08%  "
$ws.Range("E41").Value = "  +3.30%  "
$ws.Range("D42").Value = "14.52"
$ws.Range("E42").Value = "  +3.82%  "
$ws.Range("E43").Value = "  +3.71%  "
$ws.Range("E44").Value = "  -1.79%  "
$ws.Range("E45").Value = "  +3.37%  "
$ws.Range("E46").Value = "  +3.81%  "
$ws.Range("D47").Value = "6.31"
$ws.Range("E47").Value = "  +3.71%  "
$ws.Range("D48").Value = "2.029.05"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "1.01"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "3.40"
$ws.Range("E50").Value = "  +15.70%  "
$ws.Range("D51").Value = "0.0668"
$ws.Range("E51").Value = "  -0.15%  "

Write-Output "Applied cryptos update"
